$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.970.69'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.363.24'
$ws.Range('E3').Value = '  -2.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.57'
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.81'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.98'
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.417'
$ws.Range('E11').Value = '  +2.27%  '
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.98'
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('E16').Value = '  -1.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.039.51'
$ws.Range('E17').Value = '  -2.70%  '
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.50'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.38'
$ws.Range('E21').Value = '  -3.11%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.41'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.562'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('E26').Value = '  -6.73%  '
$ws.Range('E27').Value = '  -3.51%  '
$ws.Range('E28').Value = '  -2.64%  '
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.70'
$ws.Range('E32').Value = '  -3.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.87'
$ws.Range('E33').Value = '  -1.59%  '
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('E35').Value = '  +0.83%  '
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.55'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.80'
$ws.Range('E38').Value = '  -2.78%  '
$ws.Range('E39').Value = '  -8.91%  '
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('E42').Value = '  -3.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.29'
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('E44').Value = '  -2.92%  '
$ws.Range('E45').Value = '  -5.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.491.55'
$ws.Range('E46').Value = '  -2.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.68'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.815'
$ws.Range('E51').Value = '  +0.11%  '
